$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Shikuku Emmanuel"
$ws.Range("B14").Value = "Nabwana"
$ws.Range("C14").Value = "CALATECH"
$ws.Range("D14").Value = "agile"

# The phone number ("0757003013") has a leading zero and must stay text -
# assigning it directly through .Value would make Excel coerce it to a
# number and drop the leading zero. Copy it from an existing cell in the
# same column that already holds the same value as text, which preserves
# both the text type and the existing (default) cell style.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E14").PasteSpecial() | Out-Null

$ws.Range("F14").Value = "enshikuku@gmail.com"
$ws.Range("G14").Value = "gfgfg"
